$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "total"
$ws.Range("E2").Formula = "=SUM(B2:D2)"
$ws.Range("E3:E4").Formula = "=SUM(B3:D3)"

$ws.Range("E2").Select()
